$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain number (e.g. "191.57"): Excel would
# otherwise auto-convert the assigned string into a Number cell. Briefly force
# Text format so the value commits as a string (matching the source inlineStr
# cells), then restore General formatting.
$textCells = @("D5", "D6", "D11", "D14", "D18", "D22", "D26", "D28", "D31", "D33", "D36", "D37", "D39", "D42", "D44", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.539.12"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").Value = "3.377.57"
$ws.Range("E3").Value = "  +4.50%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "191.57"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "592.85"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").Value = "0.420"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "3.959.91"
$ws.Range("E12").Value = "  +4.66%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "28.64"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "69.555.94"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "3.405.51"
$ws.Range("E17").Value = "  +5.55%  "
$ws.Range("D18").Value = "448.35"
$ws.Range("E18").Value = "  +13.10%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "74.70"
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "3.518.13"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("D26").Value = "0.520"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "23.35"
$ws.Range("E31").Value = "  +3.41%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  +5.08%  "
$ws.Range("D37").Value = "165.47"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").Value = "27.32"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "6.51"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "2.743.71"
$ws.Range("E43").Value = "  +5.57%  "
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "342.90"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("D48").Value = "40.77"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "0.0285"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "32.98"
$ws.Range("E50").Value = "  +7.95%  "
$ws.Range("E51").Value = "  +5.43%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}
